# edit.ps1
# Applies three changes described by the commit diff:
#   1. Split/replace the poster title text into three runs with new wording.
#   2. Regenerate the wp14:editId on the inline picture's <wp:inline> element.
#   3. Simplify the footer paragraph ("So" + proofErr markers + " G. Presented...")
#      into a single clean run with no proofErr markup.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: "Title: Initial Results from a Phase 1 Study of
#    Cannabidiol and Tacrolimus in Healthy Subjects"
#    -> three runs: "Title: " | "Initial " | "results – Cannabidiol increased
#    exposure of tacrolimus in healthy subjects"
# ---------------------------------------------------------------------------
$titleFind = $d.Content
$foundTitle = $titleFind.Find.Execute(
    "Title: Initial Results from a Phase 1 Study of Cannabidiol and Tacrolimus in Healthy Subjects",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundTitle) {
    # Build a brand-new Range from the Find hit's coordinates -- re-using the
    # Find range object directly with InsertXML produces bad results in this
    # host, so we re-anchor via $d.Range(start, end) first.
    $titleRange = $d.Range($titleFind.Start, $titleFind.End)

    # The en dash (U+2013) is built via [char] to avoid any source-encoding
    # ambiguity in the literal script text.
    $titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/></w:rPr>' +
                '<w:t xml:space="preserve">Title: </w:t></w:r>' +
                '<w:r><w:t xml:space="preserve">Initial </w:t></w:r>' +
                '<w:r><w:t>results ' + [char]0x2013 + ' Cannabidiol increased exposure of tacrolimus in healthy subjects</w:t></w:r>' +
                '</w:p>'

    $titleRange.InsertXML($titleXml) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Inline picture: regenerate wp14:editId from 54F41A41 to 6F09A39F while
#    leaving everything else (anchorId, image relationship, sizing, alt text)
#    untouched. Located dynamically via InlineShapes(1) / the picture's docPr
#    name so we don't depend on a hard-coded paragraph index.
# ---------------------------------------------------------------------------
if ($d.InlineShapes.Count -ge 1) {
    $shape = $d.InlineShapes(1)
    $shapeParaRange = $shape.Range.Paragraphs(1).Range

    $fullXml = $d.Content.XML()
    $markerIdx = $fullXml.IndexOf('name="Picture 3"')
    if ($markerIdx -ge 0) {
        $prefix = $fullXml.Substring(0, $markerIdx)
        $pStart = $prefix.LastIndexOf("<w:p ")
        $pEndTag = $fullXml.IndexOf("</w:p>", $markerIdx)
        $pEnd = $pEndTag + 6

        $paraXml = $fullXml.Substring($pStart, $pEnd - $pStart)
        $paraXmlNew = $paraXml.Replace('wp14:editId="54F41A41"', 'wp14:editId="6F09A39F"')

        if ($paraXmlNew -ne $paraXml) {
            $shapeParaRange.InsertXML($paraXmlNew) | Out-Null
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Footer: drop the gramStart/gramEnd <w:proofErr/> markers and the run
#    split around "So" / " G. Presented..." -- collapse to a single run.
#    The host's Range.XML() accessor already normalizes away proofErr marks
#    and merges same-format adjacent runs, so we round-trip the footer's own
#    XML back into itself to pick up that normalized (target) shape while
#    preserving every untouched attribute (paraId, rsidR, pPr, ...).
# ---------------------------------------------------------------------------
for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $section = $d.Sections($si)
    $footer = $section.Footers(1)
    if ($footer.Exists) {
        $footerText = $footer.Range.Text
        if ($footerText -like "*So G. Presented at ASCPT Annual Meeting 2024, March 27, 2024.*") {
            $footerFullXml = $footer.Range.XML()
            $markerIdx = $footerFullXml.IndexOf("So G. Presented")
            if ($markerIdx -ge 0) {
                $prefix = $footerFullXml.Substring(0, $markerIdx)
                $pStart = $prefix.LastIndexOf("<w:p ")
                $pEndTag = $footerFullXml.IndexOf("</w:p>", $markerIdx)
                $pEnd = $pEndTag + 6

                $footerParaXml = $footerFullXml.Substring($pStart, $pEnd - $pStart)

                $footerRange = $footer.Range.Paragraphs(1).Range
                $footerRange.InsertXML($footerParaXml) | Out-Null
            }
        }
    }
}

Write-Output "Done."
